$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 9 (B9) which already has the desired style, into B7/B8
$ws.Range("A9:H9").Copy() | Out-Null
$ws.Range("A7:H7").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("A9:H9").Copy() | Out-Null
$ws.Range("A8:H8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("B7").Value = "23"
$ws.Range("B8").Value = "64"

$ws.Range("D7").Select() | Out-Null
